# Update "想去人数" (column F) counts for a handful of events across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 7049
$ws1.Range("F19").Value = 1011
$ws1.Range("F21").Value = 291
$ws1.Range("F28").Value = 2247

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 72
$ws2.Range("F4").Value = 57

# Sheet "全部类型" (All types) - aggregated view, mirrors the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 7049
$ws4.Range("F20").Value = 1011
$ws4.Range("F22").Value = 291
$ws4.Range("F25").Value = 72
$ws4.Range("F27").Value = 57
$ws4.Range("F31").Value = 2247
